# Weekly update: a new price observation (week of 2021-10-02, serial 44476)
# was inserted into the daily log right before the existing row for
# 2021-08-03 (serial 44384), shifting every subsequent record down by one
# row. The sheet's used range grows from A1:R112 to A1:R113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 98 (pushes old rows 98..112 down to 99..113,
# and mirrors the formatting - incl. the date number format on column D -
# of the row above, same as Excel's native "Insert Copied Cells/Rows"
# behaviour).
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A98").Value = 6
$ws.Range("B98").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C98").Value = "Metropolitana"
$ws.Range("D98").Value = 44476
$ws.Range("E98").Value = 13
$ws.Range("F98").Value = 100112029
$ws.Range("G98").Value = "Orégano"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 34
$ws.Range("K98").Value = 8500
$ws.Range("L98").Value = 9000
$ws.Range("M98").Value = 8735
$ws.Range("N98").Value = "`$/docena de atados"
$ws.Range("O98").Value = "Región Metropolitana"
$ws.Range("P98").Value = 2912
$ws.Range("Q98").Value = 3
$ws.Range("R98").Value = "Hortaliza"
